$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lũy kế tháng SÓC TRĂNG")

# Update last_edited_time text (shared string was used by D2,D3,D6,D8,D11,D13,
# all of which held the same "2024-07-28T16:31:00.000Z" timestamp)
$ws.Range("D2").Value = "2024-07-31T18:24:00.000Z"
$ws.Range("D3").Value = "2024-07-31T18:24:00.000Z"
$ws.Range("D6").Value = "2024-07-31T18:24:00.000Z"
$ws.Range("D8").Value = "2024-07-31T18:24:00.000Z"
$ws.Range("D11").Value = "2024-07-31T18:24:00.000Z"
$ws.Range("D13").Value = "2024-07-31T18:24:00.000Z"

# Update numeric figures on row 13
$ws.Range("S13").Value = 183660000
$ws.Range("W13").Value = 146528000
$ws.Range("AA13").Value = 23000000
$ws.Range("AE13").Value = 330188000
$ws.Range("AH13").Value = 293188000
$ws.Range("AK13").Value = 30
$ws.Range("AN13").Value = 37000000
$ws.Range("AQ13").Value = 316188000
